# Update forecast values on the "Forecast Comparison" sheet (Removed Auto Arima)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Row => @(D, E, F, G)
$data = @{
    2  = @(12, 15, 19, 25)
    3  = @(12, 15, 20, 30)
    4  = @(11, 14, 18, 25)
    5  = @(12, 15, 20, 29)
    6  = @(11, 14, 19, 29)
    7  = @(11, 14, 19, 29)
    8  = @(12, 14, 20, 31)
    9  = @(12, 14, 20, 31)
    10 = @(12, 14, 19, 30)
    11 = @(12, 14, 20, 30)
    12 = @(12, 14, 20, 31)
    13 = @(13, 15, 22, 34)
    14 = @(12, 13, 20, 31)
    15 = @(11, 12, 19, 30)
    16 = @(11, 13, 19, 30)
    17 = @(11, 12, 18, 29)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("D$row").Value = $values[0]
    $ws.Range("E$row").Value = $values[1]
    $ws.Range("F$row").Value = $values[2]
    $ws.Range("G$row").Value = $values[3]
}
